# edit.ps1 - applies the two changes described by the commit "9 and 10 added":
#   1. Refresh the cached "datetimeFigureOut" date text (20-09-2025 -> 05-10-2025)
#      on the Slide Master, all 11 Slide Layouts and the Notes Master.
#   2. Nudge the picture (shape id 4, "Picture 3") on slide 11 (sldId 266,
#      creationId {B8641008-E11F-2E46-EF75-9868EA0C5A59}) down slightly:
#      Top 2129913 EMU -> 2148201 EMU (Left/size unchanged).

$p = $ppt.ActivePresentation

$oldDate = "20-09-2025"
$newDate = "05-10-2025"
$ppPlaceholderDate = 16

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $ok = $false
            try {
                if ($sh.PlaceholderFormat.Type -eq $ppPlaceholderDate) { $ok = $true }
            } catch { $ok = $false }
            if ($ok) {
                # Always (re)cache the auto date/time field text to the new
                # value - some layouts report a stale/unrelated cached read
                # for this field, so don't gate the write on the read value.
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# 1a. Slide Master date placeholder.
Update-DatePlaceholder $p.SlideMaster.Shapes

# 1b. Every Slide Layout's date placeholder.
$layouts = $p.SlideMaster.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $lay = $layouts.Item($li)
    Update-DatePlaceholder $lay.Shapes
}

# 1c. Notes Master date placeholder.
Update-DatePlaceholder $p.NotesMaster.Shapes

# 2. Move the picture on slide 11 down by 18288 EMU (2129913 -> 2148201).
$s11 = $p.Slides.Item(11)
for ($i = 1; $i -le $s11.Shapes.Count; $i++) {
    $sh = $s11.Shapes.Item($i)
    if ($sh.Id -eq 4 -and $sh.Type -eq 13) {
        $sh.Top = 2148201 / 12700
    }
}
